$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 100 -> 0M
$t.Cell(1,1).Range.Text = "0M"

# Row 2: 0.02 -> 0M
$t.Cell(2,1).Range.Text = "0M"

# Row 3: 367 -> 0M
$t.Cell(3,1).Range.Text = "0M"

# Row 4: 3 -> 518
$t.Cell(4,1).Range.Text = "518"

# Row 5: 0.00003 -> 0.00002
$t.Cell(5,1).Range.Text = "0.00002"

# Row 6: 0.00003 -> 0.00006
$t.Cell(6,1).Range.Text = "0.00006"

# Row 9: 0.00003 -> 0.00004
$t.Cell(9,1).Range.Text = "0.00004"

# Row 10: 0.00003 -> 0.00004
$t.Cell(10,1).Range.Text = "0.00004"

# Row 11: 0.00003 -> 0.00004
$t.Cell(11,1).Range.Text = "0.00004"

# Row 12: 0.00008 -> 0.01564
$t.Cell(12,1).Range.Text = "0.01564"

# Row 44: collapse multi-run tab-separated text down to "100"
$t.Cell(44,1).Range.Text = "100"

# Row 45: collapse multi-run tab-separated text down to "0.02"
$t.Cell(45,1).Range.Text = "0.02"

# Row 46: collapse multi-run tab-separated text down to "367"
$t.Cell(46,1).Range.Text = "367"
